# Fix header labels on existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Clone the existing header style (bold + border) and date-column style
# from "Weekly Quantity" so the new sheet reuses the same style entries.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:D1"))
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A31"))

# Header row labels
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = New-Object 'object[,]' 30,4
$data[0,0] = 45193.99999999999; $data[0,1] = 64; $data[0,2] = -96.0659480322405; $data[0,3] = 229.9307406276706
$data[1,0] = 45207.99999999999; $data[1,1] = 70; $data[1,2] = -64.92108787613033; $data[1,3] = 231.5951540091356
$data[2,0] = 45214.99999999999; $data[2,1] = 73; $data[2,2] = -88.13942659371352; $data[2,3] = 243.1788897195244
$data[3,0] = 45228.99999999999; $data[3,1] = 78; $data[3,2] = -70.70305724727758; $data[3,3] = 235.7134963446877
$data[4,0] = 45235.99999999999; $data[4,1] = 81; $data[4,2] = -59.91895817435565; $data[4,3] = 239.5519565755657
$data[5,0] = 45242.99999999999; $data[5,1] = 84; $data[5,2] = -78.44897531073381; $data[5,3] = 240.161510200797
$data[6,0] = 45249.99999999999; $data[6,1] = 87; $data[6,2] = -67.765475704302; $data[6,3] = 232.5326713558642
$data[7,0] = 45256.99999999999; $data[7,1] = 90; $data[7,2] = -75.4716185098495; $data[7,3] = 248.1236578030863
$data[8,0] = 45263.99999999999; $data[8,1] = 93; $data[8,2] = -64.53853042203635; $data[8,3] = 249.5548404848613
$data[9,0] = 45270.99999999999; $data[9,1] = 96; $data[9,2] = -58.27092239456469; $data[9,3] = 251.5501032903905
$data[10,0] = 45277.99999999999; $data[10,1] = 99; $data[10,2] = -64.64097952451326; $data[10,3] = 255.217624736096
$data[11,0] = 45298.99999999999; $data[11,1] = 107; $data[11,2] = -60.69735333354729; $data[11,3] = 258.4652017386988
$data[12,0] = 45305.99999999999; $data[12,1] = 110; $data[12,2] = -50.35430749492587; $data[12,3] = 269.8888841869829
$data[13,0] = 45319.99999999999; $data[13,1] = 116; $data[13,2] = -40.31755873628241; $data[13,3] = 271.3064869648707
$data[14,0] = 45326.99999999999; $data[14,1] = 119; $data[14,2] = -41.48358617510715; $data[14,3] = 277.4916841213555
$data[15,0] = 45333.99999999999; $data[15,1] = 122; $data[15,2] = -46.55593724514933; $data[15,3] = 278.4491714302915
$data[16,0] = 45340.99999999999; $data[16,1] = 124; $data[16,2] = -40.27970292832994; $data[16,3] = 271.4397844652273
$data[17,0] = 45347.99999999999; $data[17,1] = 127; $data[17,2] = -23.80085336501163; $data[17,3] = 273.884010905565
$data[18,0] = 45354.99999999999; $data[18,1] = 130; $data[18,2] = -21.41932295122118; $data[18,3] = 280.8739992605001
$data[19,0] = 45361.99999999999; $data[19,1] = 133; $data[19,2] = -26.10904850297245; $data[19,3] = 293.1279027323242
$data[20,0] = 45368.99999999999; $data[20,1] = 136; $data[20,2] = -16.64751249526828; $data[20,3] = 286.8376243235055
$data[21,0] = 45375.99999999999; $data[21,1] = 139; $data[21,2] = -8.586365764754749; $data[21,3] = 274.6735401242062
$data[22,0] = 45382.99999999999; $data[22,1] = 142; $data[22,2] = -22.83238959761452; $data[22,3] = 304.7508092632249
$data[23,0] = 45389.99999999999; $data[23,1] = 144; $data[23,2] = -13.61536169940072; $data[23,3] = 295.3988581391383
$data[24,0] = 45396.99999999999; $data[24,1] = 147; $data[24,2] = -8.861155452744383; $data[24,3] = 299.9306198846508
$data[25,0] = 45403.99999999999; $data[25,1] = 150; $data[25,2] = -14.44785939633723; $data[25,3] = 314.8972144738126
$data[26,0] = 45410.99999999999; $data[26,1] = 153; $data[26,2] = -7.494539147398013; $data[26,3] = 305.572135635525
$data[27,0] = 45417.99999999999; $data[27,1] = 156; $data[27,2] = -4.093376903688263; $data[27,3] = 324.486088085379
$data[28,0] = 45424.99999999999; $data[28,1] = 159; $data[28,2] = 4.526512648885775; $data[28,3] = 297.0156743454094
$data[29,0] = 45431.99999999999; $data[29,1] = 162; $data[29,2] = 10.51746002138638; $data[29,3] = 321.2947088869373

$wsForecast.Range("A2:D31").Value = $data
